$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("E1").Value = "Target_Mean"
$ws.Range("F1").Value = "Target_Std"
$ws.Range("G1").Value = "Target_Min"
$ws.Range("H1").Value = "Target_Max"

# Copy the header style (s="1") from an existing header cell onto the new header cells
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Helper ranges that hold text that looks numeric ("59.9%", "0.38", ...).
# Force them to be stored as text (not auto-converted to numbers/percentages)
# by temporarily switching to a text number format, then clearing the
# formatting again so the cell keeps its original (default) style.
$textRanges = @("C2:C4", "E2:F4")
foreach ($rng in $textRanges) {
    $ws.Range($rng).NumberFormat = "@"
}

# --- Row 2 (Train) ---
$ws.Range("B2").Value = 534
$ws.Range("C2").Value = "59.9%"
$ws.Range("D2").Value = 532
$ws.Range("E2").Value = "0.38"
$ws.Range("F2").Value = "0.49"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1

# --- Row 3 (Val) ---
$ws.Range("B3").Value = 178
$ws.Range("C3").Value = "20.0%"
$ws.Range("D3").Value = 169
$ws.Range("E3").Value = "0.41"
$ws.Range("F3").Value = "0.49"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1

# --- Row 4 (Test) ---
$ws.Range("B4").Value = 179
$ws.Range("D4").Value = 165
$ws.Range("E4").Value = "0.37"
$ws.Range("F4").Value = "0.48"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1

# Remove the temporary text number format again so these cells end up with
# the default (unstyled) look, matching the rest of the data rows.
foreach ($rng in $textRanges) {
    $ws.Range($rng).ClearFormats()
}
